$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the cells keep their text (string) nature, since the original
# values were stored as shared strings (text), not numbers/dates.
$ws.Range("A2:I2").NumberFormat = "@"

$ws.Range("A2").Value = "Victor Manuel"
$ws.Range("B2").Value = "Bonilla Gutierrez"
$ws.Range("C2").Value = "Permiso por Protección Temporal"
$ws.Range("D2").Value = "4073477"
$ws.Range("E2").Value = "ADSO"
$ws.Range("F2").Value = "Tecnologo"
$ws.Range("G2").Value = "Desarrollo Web"
$ws.Range("H2").Value = "03:34 p." + [char]0x202F + "m."
$ws.Range("I2").Value = "A tiempo"
